$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# SupIm sheet: extend the single "year 1" row of capacity-factor data into a
# full 12-row timeseries (rows 2..13 -> years 1..12), keeping the formatting
# of the template row (row 3) on every new row.
# ---------------------------------------------------------------------------
$supim = $wb.Worksheets.Item("SupIm")

for ($i = 4; $i -le 14; $i++) {
    $year = $i - 2
    $supim.Cells.Item($i, 1).Value = $year
    $supim.Cells.Item($i, 2).Value = 0.481
    $supim.Cells.Item($i, 3).Value = 0.3
    $supim.Cells.Item($i, 4).Value = 0.207

    $supim.Range("A3:D3").Copy()
    $supim.Range("A" + $i + ":D" + $i).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

[void]$supim.Range("H16").Select()

# ---------------------------------------------------------------------------
# Demand sheet: replace the single annual demand value with a flat monthly
# timeseries (12 equal instalments of the previous yearly total).
# ---------------------------------------------------------------------------
$demand = $wb.Worksheets.Item("Demand")

$demand.Cells.Item(3, 2).Value = 565916667

for ($i = 4; $i -le 14; $i++) {
    $year = $i - 2
    $demand.Cells.Item($i, 1).Value = $year
    $demand.Cells.Item($i, 2).Value = 565916667

    $demand.Range("B3").Copy()
    $demand.Range("B" + $i).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

[void]$demand.Range("D13").Select()
$demand.Activate()
